$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

function Get-ParaByExactText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($text + "`r")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. Title change (Heading1) and the identical bold meta-title further down.
# ---------------------------------------------------------------------
Replace-AllText "Play Ghostbusters Plus Slot for Free - Exciting Gameplay" "Play Ghostbuster Plus for Free"

# ---------------------------------------------------------------------
# 2. "What we like" bullet list changes.
# ---------------------------------------------------------------------

# Remove the "Familiar characters and symbols..." bullet entirely.
$p = Get-ParaByExactText "Familiar characters and symbols from the Ghostbusters franchise"
if ($p -ne $null) {
    $p.Range.Delete()
}

# "Random events during gameplay for added excitement" -> "Fun and varied gameplay with bonus features"
Replace-AllText "Random events during gameplay for added excitement" "Fun and varied gameplay with bonus features"

# "Fun and varied gameplay experience" -> "Appeal to fans of the Ghostbusters film franchise"
Replace-AllText "Fun and varied gameplay experience" "Appeal to fans of the Ghostbusters film franchise"

# New bullet "Intuitive gameplay mechanics" right after the "Appeal to fans..." bullet.
$p = Get-ParaByExactText "Appeal to fans of the Ghostbusters film franchise"
if ($p -ne $null) {
    $srcFormatted = $p.Range.FormattedText
    $p.Range.InsertParagraphAfter()
    $idx = $p.Index
    $newPara = $d.Paragraphs($idx + 1)
    $newPara.Range.FormattedText = $srcFormatted
    $newPara.Range.Find.Execute(
        "Appeal to fans of the Ghostbusters film franchise", $true, $false, $false, $false, $false,
        $true, 1, $false, "Intuitive gameplay mechanics", 2
    ) | Out-Null
}

# ---------------------------------------------------------------------
# 3. "What we don't like" changes.
# ---------------------------------------------------------------------
Replace-AllText "Not significantly different from previous Ghostbusters iterations" "Not significantly different from previous iterations"
Replace-AllText "Competition from other cinema-themed slots" "Competition from other cinema-themed slot games"

# ---------------------------------------------------------------------
# 4. Meta description (italic paragraph) change.
# ---------------------------------------------------------------------
Replace-AllText "Read our review of Ghostbusters Plus and play this exciting slot game for free. Fun gameplay and detailed graphics await fans of the franchise." "Read our review of the Ghostbuster Plus slot game and play for free. Enjoy improved graphics and fun gameplay."
